$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the teacher / bill details (Label wise bill generate - Individual)
$ws.Range("A3").Value = "নাম: Dr. S. M. Rabiul Alam (Hum)"
$ws.Range("A4").Value = "পদবী: অধ্যাপক"
$ws.Range("F5").Value = "বিভাগ :হুম"

# Invigilation row: 1 day/unit -> amount is computed by existing formula (K26*G26)
$ws.Range("G26").Value = 1

# Amount in words for the grand total
$ws.Range("A32").Value = "কথায়:দুই হাজার সাতশো টাকা মাত্র।"

# Move the active selection to B5 as in the saved file
$ws.Range("B5").Select() | Out-Null
